# Update examples and documentation
# Travis_County_2017_bg_SVI.xlsx edit: reorder factor-2 loading rows by
# loading magnitude and refresh the recomputed numeric values, and update
# the example variable-list strings on "Significant Components" and
# "Included and Excluded" to match the new ordering.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Significant Components": update example variable lists
# ---------------------------------------------------------------
$wsSig = $wb.Worksheets.Item("Significant Components")

$wsSig.Range("C2").Value = "['PPUNIT' 'QNOHLTH' 'QSERV' 'QEXTRCT' 'QESL' 'QHISPC' 'QEDLESHI' 'QFHH'`n 'PERCAP']"
$wsSig.Range("C4").Value = "['MEDAGE' 'QAGEDEP' 'QSSBEN']"
$wsSig.Range("C5").Value = "['QAGEDEP' 'QFEMALE' 'QFEMLBR']"

# ---------------------------------------------------------------
# Sheet "Loading Factors": rows 2-20 re-sorted by the new factor
# ordering, with refreshed loading values for each variable.
# ---------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Loading Factors")

$loadData = @(
@{Row=2; Label="PPUNIT"; B=0.7302277669272809; C=-0.004620269588076043; D=-0.1512975051829958; E=0.05458729952269146; F=-0.4752714268289649},
    @{Row=3; Label="QNOHLTH"; B=0.6889886602202583; C=0.4149592268437005; D=-0.1190896980125687; E=-0.1154839672367; F=0.2786568980252105},
    @{Row=4; Label="QSERV"; B=0.5817049374264434; C=0.3577287403200355; D=-0.2244058813072648; E=-0.03201466882836596; F=0.2754259478856783},
    @{Row=5; Label="QEXTRCT"; B=0.7677511934101074; C=0.1449653485907693; D=0.01129678766226058; E=-0.2382017664847616; F=0.09278788606079744},
    @{Row=6; Label="QESL"; B=0.8009661392984533; C=0.151795080540303; D=-0.03424774568502074; E=-0.2374898269416714; F=0.2037776736638281},
    @{Row=7; Label="QHISPC"; B=0.8328587039996096; C=0.3339037448970698; D=-0.1364957895465116; E=-0.1269704316134753; F=0.09670565299798836},
    @{Row=8; Label="QEDLESHI"; B=0.8777939322771465; C=0.2130779229266954; D=-0.01839360962630974; E=-0.1076883447776996; F=0.1846338469805828},
    @{Row=9; Label="QFHH"; B=0.5630560154859281; C=0.3008187133860031; D=-0.09551317270938871; E=0.2634165278540286; F=-0.03178269088452637},
    @{Row=10; Label="PERCAP"; B=0.4895374274546244; C=0.7214540806533332; D=-0.2685725744248457; E=0.05482089686608119; F=0.183125353910591},
    @{Row=11; Label="QRICH"; B=0.2150971475063948; C=0.8701305525241234; D=-0.1729034330861709; E=-0.01425153555799462; F=0.2948729863855478},
    @{Row=12; Label="MDHSEVAL"; B=0.3857664787720542; C=0.8013369875455929; D=-0.03601197820671726; E=-0.02873185417247733; F=-0.03016150210197751},
    @{Row=13; Label="MEDAGE"; B=-0.3105472641016062; C=-0.2465384264432842; D=0.7910361932243746; E=-0.0128947456807908; F=-0.2717085837049965},
    @{Row=14; Label="QAGEDEP"; B=-0.04282432074929374; C=-0.1184872976686267; D=0.6543380974999534; E=0.6427560499114501; F=-0.1139548224937056},
    @{Row=15; Label="QSSBEN"; B=0.01836498138439658; C=-0.05367045618418921; D=0.7773306998288915; E=0.1362099740686128; F=-0.1455945116013038},
    @{Row=16; Label="QFEMALE"; B=-0.04556235957520197; C=-0.05659266691271623; D=0.1671554291910123; E=0.877843202586366; F=-0.02420501690989349},
    @{Row=17; Label="QFEMLBR"; B=-0.2416846119055139; C=0.08178103489040496; D=-0.0295955310337048; E=0.7849928999651834; F=0.003645846047376428},
    @{Row=18; Label="QRENTER"; B=0.01759023509580655; C=0.2288188122079551; D=-0.423414426863078; E=-0.09662324553205429; F=0.7659446309671079},
    @{Row=19; Label="QNOAUTO"; B=0.1660837603234652; C=0.06295541717908909; D=-0.1064187385486532; E=-0.01530326288287446; F=0.6312035337968215},
    @{Row=20; Label="QPOVTY"; B=0.3701334413584173; C=0.1576334738980672; D=-0.3817953695838431; E=0.08010852895572017; F=0.4611334674601166}
)

foreach ($item in $loadData) {
    $wsLoad.Cells.Item($item.Row, 1).Value = $item.Label
    $wsLoad.Cells.Item($item.Row, 2).Value = $item.B
    $wsLoad.Cells.Item($item.Row, 3).Value = $item.C
    $wsLoad.Cells.Item($item.Row, 4).Value = $item.D
    $wsLoad.Cells.Item($item.Row, 5).Value = $item.E
    $wsLoad.Cells.Item($item.Row, 6).Value = $item.F
}

# ---------------------------------------------------------------
# Sheet "All Refactor Variances": refreshed SS loadings / variance
# figures (rows 2-5, columns B-R).
# ---------------------------------------------------------------
$wsAllVar = $wb.Worksheets.Item("All Refactor Variances")

$allVarData = @(
@{Row=2; B=5.346593380343492; C=2.352336411856979; D=2.340013488429634; E=2.207229445501692; F=2.005491058151331; G=1.966280396698093; H=1.011763029505047; I=4.839746430455178; J=3.42202125223686; K=2.232417339796182; L=2.056337895353981; M=2.044708469958762; N=5.117234762546761; O=2.680150570577836; P=2.243476579107345; Q=2.066288182496249; R=1.90171479706631},
    @{Row=3; B=0.198021977049759; C=0.08712357080951776; D=0.08666716623813461; E=0.0817492387222849; F=0.07427744659819745; G=0.07282519987770715; H=0.03747270479648322; I=0.2304641157359609; J=0.16295339296366; K=0.106305587609342; L=0.0979208521597134; M=0.09736706999803627; N=0.2693281453971979; O=0.1410605563462019; P=0.1180777146898602; Q=0.1087520096050658; R=0.1000902524771742},
    @{Row=4; B=0.198021977049759; C=0.2851455478592767; D=0.3718127140974113; E=0.4535619528196962; F=0.5278393994178937; G=0.6006645992956009; H=0.6381373040920841; I=0.2304641157359609; J=0.3934175086996209; K=0.4997230963089629; L=0.5976439484686763; M=0.6950110184667125; N=0.2693281453971979; O=0.4103887017433998; P=0.5284664164332601; Q=0.6372184260383258; R=0.7373086785155},
    @{Row=5; B=0.3103124919667509; C=0.1365279388163549; D=0.135812725070886; E=0.1281060332910554; F=0.1163972802120954; G=0.1141215212003941; H=0.05872200944246296; I=0.3315977871032831; J=0.2344615964839767; K=0.152955255074756; L=0.1408910787856859; M=0.1400942825522983; N=0.3652854676001701; O=0.1913181825422342; P=0.1601469209986764; Q=0.1474986159447186; R=0.1357508129142007}
)

$allVarCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

foreach ($item in $allVarData) {
    foreach ($col in $allVarCols) {
        $wsAllVar.Range("$col$($item.Row)").Value = $item[$col]
    }
}

# ---------------------------------------------------------------
# Sheet "Final Variances": refreshed SS loadings / variance figures
# (rows 2-5, columns B-F) mirroring the factor-2 block above.
# ---------------------------------------------------------------
$wsFinalVar = $wb.Worksheets.Item("Final Variances")

$finalVarData = @(
@{Row=2; B=5.117234762546761; C=2.680150570577836; D=2.243476579107345; E=2.066288182496249; F=1.90171479706631},
    @{Row=3; B=0.2693281453971979; C=0.1410605563462019; D=0.1180777146898602; E=0.1087520096050658; F=0.1000902524771742},
    @{Row=4; B=0.2693281453971979; C=0.4103887017433998; D=0.5284664164332601; E=0.6372184260383258; F=0.7373086785155},
    @{Row=5; B=0.3652854676001701; C=0.1913181825422342; D=0.1601469209986764; E=0.1474986159447186; F=0.1357508129142007}
)

$finalVarCols = @("B","C","D","E","F")

foreach ($item in $finalVarData) {
    foreach ($col in $finalVarCols) {
        $wsFinalVar.Range("$col$($item.Row)").Value = $item[$col]
    }
}

# ---------------------------------------------------------------
# Sheet "Included and Excluded": update the combined "include" list
# to the new variable ordering.
# ---------------------------------------------------------------
$wsInclExcl = $wb.Worksheets.Item("Included and Excluded")

$wsInclExcl.Range("B2").Value = "[['PPUNIT', 'QNOHLTH', 'QSERV', 'QEXTRCT', 'QESL', 'QHISPC', 'QEDLESHI', 'QFHH', 'PERCAP', 'QRICH', 'MDHSEVAL', 'MEDAGE', 'QAGEDEP', 'QSSBEN', 'QFEMALE', 'QFEMLBR', 'QRENTER', 'QNOAUTO', 'QPOVTY']]"
